# The "Enterprises density (per 1000 people)" / "8.5" row needs to appear
# before the "Enterprises (absolute #)" / "253080" row, i.e. swap the
# content of rows 10 and 11 on the Summary sheet:
#   row 10: "Enterprises (absolute #)" / "253080"  ->  "Enterprises density (per 1000 people)" / "8.5"
#   row 11: "Enterprises density (per 1000 people)" / "8.5"  ->  "Enterprises (absolute #)" / "253080"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking values to stay stored as text (as they were in
# the original file: cells D10/D11 are t="s", not real numbers).
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"

$ws.Range("A10").Value = "Enterprises density (per 1000 people)"
$ws.Range("D10").Value = "8.5"

$ws.Range("A11").Value = "Enterprises (absolute #)"
$ws.Range("D11").Value = "253080"
